# Generate Report for Handoff
# Replaces the two handed-back source files (2db4db43..., dccdc1f8...) with a
# single newly-handed-off file (b488e0a4...) and a freshly queued file
# (ffffb1983d0c...), updates status text/timestamps, and clears the
# now-stale "Latest Target File" / "Latest Handback File" columns on the
# per-language sheets.

$wb = $excel.ActiveWorkbook

$oldFile1 = "2db4db43-ed49-4db6-94a7-647e8cb93e42.md"
$oldFile2 = "dccdc1f8-1fc0-4ee0-ac3a-7fe9156b3bc5.md"
$newFile1 = "b488e0a4-e7d8-4343-a46b-53f4ea708df8.md"
$newFile2 = "ffffb1983d0c-e3f9-4ca6-80eb-d06f1cff7aff.md"

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

$newHoDate = "2016-08-30 11:12:19"
$newHandoffDate = "2016-08-30 11:12:14"
$zeroDate = "0001-01-01 00:00:00"

$newXlfZh = "b488e0a4-e7d8-4343-a46b-53f4ea708df8.fa112066d26f78bbcbadb52fba6ff71b07da0b4a.zh-cn.xlf"
$newXlfDe = "b488e0a4-e7d8-4343-a46b-53f4ea708df8.fa112066d26f78bbcbadb52fba6ff71b07da0b4a.de-de.xlf"

$srcRepo = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06b8b5fed0864774f689490a00885d9a7d5f693e/e2e/"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $newFile1
$ws.Range("C2").Value = ".md"
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Range("G2").Value = $newHoDate

$ws.Range("A3").Value = $newFile2
$ws.Range("C3").Value = ".md"
$ws.Range("E3").Value = $newStatus
$ws.Range("F3").Value = $newStatus
$ws.Range("G3").Value = $newHoDate

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), ($srcRepo + $newFile1), "", "", ("e2e\" + $newFile1))
$ws.Hyperlinks.Add($ws.Range("B3"), ($srcRepo + $newFile2), "", "", ("e2e\" + $newFile2))

$ws.Columns.Item(5).ColumnWidth = 16.3
$ws.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $newFile1
$ws.Range("C2").Value = $newStatus
$ws.Range("G2").Value = $newXlfZh
$ws.Range("H2").Value = $newHandoffDate
$ws.Range("I2").Hyperlinks.Delete()
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = $zeroDate

$ws.Range("A3").Value = $newFile2
$ws.Range("C3").Value = $newStatus
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = $newXlfZh
$ws.Range("H3").Value = $newHandoffDate
$ws.Range("I3").Hyperlinks.Delete()
$ws.Range("I3").Value = ""
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = $zeroDate

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), ($srcRepo + $newFile1), "", "", $newFile1)
$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), ($srcRepo + $newFile2), "", "", $newFile2)

$ws.Columns.Item(3).ColumnWidth = 16.3
$ws.Columns.Item(9).ColumnWidth = 17.8
$ws.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $newFile1
$ws.Range("C2").Value = $newStatus
$ws.Range("G2").Value = $newXlfDe
$ws.Range("H2").Value = $newHoDate
$ws.Range("I2").Hyperlinks.Delete()
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = $zeroDate

$ws.Range("A3").Value = $newFile2
$ws.Range("C3").Value = $newStatus
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = $newXlfDe
$ws.Range("H3").Value = $newHoDate
$ws.Range("I3").Hyperlinks.Delete()
$ws.Range("I3").Value = ""
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = $zeroDate

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), ($srcRepo + $newFile1), "", "", $newFile1)
$ws.Range("A3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), ($srcRepo + $newFile2), "", "", $newFile2)

$ws.Columns.Item(3).ColumnWidth = 16.3
$ws.Columns.Item(9).ColumnWidth = 17.8
$ws.Columns.Item(10).ColumnWidth = 20.8

Write-Output "done"
